# Bug fix update 2024
#
# Column A (rows 2-17) held numeric date-serial values (formatted with a
# date number format). The fix turns those into literal text labels
# ("01.01.2001", "01.04.2001", ...) instead of real dates, while keeping
# the original cell formatting untouched.
#
# Writing the text straight into .Value/.Value2/.Formula gets silently
# re-parsed back into a date serial by Excel's smart entry, so instead we
# enter each label as a text-literal formula (="01.01.2001") and then
# Copy / Paste-Special-Values it back onto itself. That converts the
# formula result into a plain literal string in place, without touching
# the cell's existing number format/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateLabels = @{
    2  = "01.01.2001"
    3  = "01.04.2001"
    4  = "01.07.2001"
    5  = "01.10.2001"
    6  = "01.01.2002"
    7  = "01.04.2002"
    8  = "01.07.2002"
    9  = "01.10.2002"
    10 = "01.01.2003"
    11 = "01.04.2003"
    12 = "01.07.2003"
    13 = "01.10.2003"
    14 = "01.01.2004"
    15 = "01.04.2004"
    16 = "01.07.2004"
    17 = "01.10.2004"
}

foreach ($row in 2..17) {
    $label = $dateLabels[$row]
    $ws.Range("A$row").Formula = "=""$label"""
}

# Freeze the formula results as literal text values in-place, preserving
# the existing (date) number format on A2:A17.
$ws.Range("A2:A17").Copy()
$ws.Range("A2:A17").PasteSpecial(-4163)
